# WIR19 tab07 - add 2019 data column (AE) to the "net sales" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2019" header cell (plain numeric literal, matching the header row's
# general style s=6 rather than the boxed header style s=8 used by B1:AD1).
$ws.Range("AE1").Value = 2019

# New 2019 data values for each country row.
$ws.Range("AE2").Value = 37
$ws.Range("AE3").Value = 79
$ws.Range("AE4").Value = 155
$ws.Range("AE5").Value = 131
$ws.Range("AE6").Value = 58

# Update the saved selection to match the author's last cursor position.
$ws.Range("Z4").Select()
